# Adds a new "2022-Q3" sheet (right after "总计") with the Q3 fund-holding
# detail data, and records the Q3 summary row at the top of the "总计"
# (totals) sheet's data table.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)      # "总计"
$q2    = $wb.Worksheets.Item(2)      # "2022-Q2" (reference sheet for layout/format)

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet right after "总计"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# Header row (matches the other quarterly sheets)
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
$cols = @("B","C","D","E","F","G","H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3.Range($cols[$i] + "1")
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Fund holding rows for 2022-Q3
$rows = @(
    @("011531","朱雀恒心一年持有期混合","59.37","93.18","5.46","3.2416",5),
    @("007493","朱雀产业臻选混合A","34.47","92.95","4.76","1.6408",5),
    @("010922","朱雀匠心一年持有期混合","15.23","91.81","9.52","1.4499",1),
    @("010141","朱雀企业优选股票A","27.92","93.51","4.90","1.3681",6),
    @("007494","朱雀产业臻选混合C","11.69","92.95","4.76","0.5564",5),
    @("008294","朱雀企业优胜股票A","11.46","93.76","4.75","0.5444",5),
    @("010142","朱雀企业优选股票C","4.96","93.51","4.90","0.2430",6),
    @("007880","朱雀产业智选混合A","3.92","93.07","4.32","0.1693",5),
    @("008295","朱雀企业优胜股票C","2.16","93.76","4.75","0.1026",5),
    @("007881","朱雀产业智选混合C","0.71","93.07","4.32","0.0307",5)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $data = $rows[$r]

    $aCell = $q3.Range("A" + $rowNum)
    $aCell.Value = $r
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    # B, C, D, E, F, G are stored as text (matches source data); H is numeric
    $q3.Range("B" + $rowNum + ":G" + $rowNum).NumberFormat = "@"
    $q3.Range("B" + $rowNum).Value = $data[0]
    $q3.Range("C" + $rowNum).Value = $data[1]
    $q3.Range("D" + $rowNum).Value = $data[2]
    $q3.Range("E" + $rowNum).Value = $data[3]
    $q3.Range("F" + $rowNum).Value = $data[4]
    $q3.Range("G" + $rowNum).Value = $data[5]
    $q3.Range("H" + $rowNum).Value = $data[6]
}

# ---------------------------------------------------------------------
# 2. Insert the 2022-Q3 summary row at the top of "总计"'s data table
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").Style = "Normal"

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 10
$total.Range("D2").Value = 9.35

# Re-number the sequential index column (A) for every data row, and make
# sure the new row's A cell carries the same look as the others.
$lastRow = $total.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $total.Range("A" + $r)
    $cell.Value = $r - 2
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

Write-Host "done"
